$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.907.25"
$ws.Range("E2").Value = "  -0.35%  "
$ws.Range("D3").Value = "3.795.84"
$ws.Range("E3").Value = "  -2.09%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'598.66"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "'168.83"
$ws.Range("E6").Value = "  -1.54%  "
$ws.Range("D7").Value = "3.795.28"
$ws.Range("E7").Value = "  -2.10%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("E10").Value = "  +1.51%  "
$ws.Range("D11").Value = "'6.52"
$ws.Range("E11").Value = "  +1.63%  "
$ws.Range("D12").Value = "'0.461"
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("D13").Value = "'0.0000274"
$ws.Range("E13").Value = "  +5.80%  "
$ws.Range("D14").Value = "'36.95"
$ws.Range("E14").Value = "  -0.29%  "
$ws.Range("D15").Value = "4.433.35"
$ws.Range("E15").Value = "  -2.04%  "
$ws.Range("D16").Value = "3.796.20"
$ws.Range("E16").Value = "  -2.07%  "
$ws.Range("D17").Value = "'19.19"
$ws.Range("E17").Value = "  +5.65%  "
$ws.Range("D18").Value = "67.871.51"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").Value = "'7.31"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("D21").Value = "'10.62"
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("D22").Value = "'467.25"
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Value = "'0.730"
$ws.Range("E23").Value = "  -1.35%  "
$ws.Range("E24").Value = "  -5.52%  "
$ws.Range("D25").Value = "'83.54"
$ws.Range("E25").Value = "  +0.38%  "
$ws.Range("E26").Value = "  +2.15%  "
$ws.Range("D27").Value = "'12.23"
$ws.Range("E27").Value = "  +1.46%  "
$ws.Range("D28").Value = "'10.35"
$ws.Range("E28").Value = "  +3.63%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").Value = "'2.94"
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("D31").Value = "3.945.79"
$ws.Range("E31").Value = "  -1.95%  "
$ws.Range("E32").Value = "  -1.62%  "
$ws.Range("E33").Value = "  -2.49%  "
$ws.Range("D34").Value = "'30.58"
$ws.Range("E34").Value = "  -2.07%  "
$ws.Range("E35").Value = "  -3.05%  "
$ws.Range("D36").Value = "3.758.04"
$ws.Range("E36").Value = "  -2.05%  "
$ws.Range("E37").Value = "  +1.01%  "
$ws.Range("D38").Value = "'0.106"
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("D39").Value = "'5.94"
$ws.Range("E39").Value = "  +0.53%  "
$ws.Range("E40").Value = "  -1.56%  "
$ws.Range("D41").Value = "'0.138"
$ws.Range("E41").Value = "  -1.62%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").Value = "'0.320"
$ws.Range("E43").Value = "  +2.16%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "'8.77"
$ws.Range("E45").Value = "  +1.68%  "
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("D47").Value = "'408.95"
$ws.Range("E47").Value = "  -3.71%  "
$ws.Range("D48").Value = "'46.29"
$ws.Range("E48").Value = "  -1.89%  "
$ws.Range("D49").Value = "'0.000279"
$ws.Range("E49").Value = "  -6.36%  "
$ws.Range("D50").Value = "'142.24"
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("D51").Value = "'0.0357"
$ws.Range("E51").Value = "  -0.01%  "
